$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "57.560.46"
$ws.Cells.Item(2, 5).Value = "  -1.70%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.434.17"
$ws.Cells.Item(3, 5).Value = "  -2.16%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.06%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "514.23"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -2.54%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "129.99"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -3.00%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.18%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -2.07%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "2.446.83"
$ws.Cells.Item(9, 5).Value = "  -1.73%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  -0.11%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0950"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -5.37%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -3.88%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  -3.57%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "2.864.02"
$ws.Cells.Item(14, 5).Value = "  -2.19%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "57.416.08"
$ws.Cells.Item(15, 5).Value = "  -1.78%  "

# Row 16
$ws.Cells.Item(16, 5).Value = "  -3.02%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  -3.44%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "2.440.19"
$ws.Cells.Item(18, 5).Value = "  -1.91%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "10.46"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -4.54%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "315.96"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -1.88%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  -2.62%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "5.68"
$ws.Cells.Item(23, 4).Style = "Normal"

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "63.31"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -1.74%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -1.59%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.998"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -0.13%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -1.69%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -3.48%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "170.14"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +2.99%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "0.0₃0723"
$ws.Cells.Item(30, 5).Value = "  -4.22%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "6.26"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -2.75%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  -2.76%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +2.28%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -0.04%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.997"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.15%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "17.70"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -3.24%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  -4.69%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "3.93"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -1.76%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "36.18"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -1.55%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "1.46"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -3.17%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.782"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -2.67%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "Bittensor"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "270.23"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -2.31%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "RenderToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "5.00"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +0.75%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "Filecoin"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "3.38"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -4.48%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.585"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -2.01%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.0908"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -0.88%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "120.63"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -5.39%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.0485"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -2.38%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  -2.95%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "16.59"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -3.94%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "1.708.55"
$ws.Cells.Item(51, 5).Value = "  -1.99%  "
